$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column BJ (column 62) with header "05-sep", copying formats from
# the previous column (BI, column 61) so the new column matches existing style.
$ws.Range("BI1:BI11").Copy()
$ws.Range("BJ1:BJ11").PasteSpecial(-4122)

$ws.Range("BJ1").Value = "05-sep"

$values = @(16, 13, 11, 13, 9, 20, 20, 17, 18, 10)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 62).Value = $values[$i]
}

$ws.Range("BK7").Select()
